# Update the "nemad" (ticker symbol) column L for rows 2-45 from the
# company name to the correct ticker symbol "شسپا".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 12).Value = "شسپا"
}

$ws.Range("L5").Select()
